$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.735.04'
$ws.Range('E2').Value = '  +1.80%  '
$ws.Range('D3').Value = '1.694.86'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Formula = "'316.75"
$ws.Range('E5').Value = '  +1.67%  '
$ws.Range('D6').Formula = "'1.0000"
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').Formula = "'0.3947"
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('D8').Formula = "'0.4056"
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('D11').Formula = "'53.25"
$ws.Range('E11').Value = '  -2.43%  '
$ws.Range('D12').Formula = "'0.08865"
$ws.Range('E12').Value = '  +1.33%  '
$ws.Range('D13').Formula = "'7.229"
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('D14').Formula = "'23.61"
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('D15').Formula = "'8.057"
$ws.Range('E15').Value = '  +8.82%  '
$ws.Range('D16').Formula = "'0.00001321"
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('D17').Value = '1.692.57'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').Formula = "'99.97"
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').Formula = "'0.07009"
$ws.Range('E19').Value = '  -0.50%  '
$ws.Range('D20').Formula = "'19.62"
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('D21').Formula = "'7.054"
$ws.Range('E21').Value = '  +4.76%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').Formula = "'14.33"
$ws.Range('E23').Value = '  +1.35%  '
$ws.Range('D24').Value = '24.731.62'
$ws.Range('E24').Value = '  +1.79%  '
$ws.Range('D25').Formula = "'3.279"
$ws.Range('E25').Value = '  +9.86%  '
$ws.Range('D26').Formula = "'2.356"
$ws.Range('E26').Value = '  +1.84%  '
$ws.Range('D28').Formula = "'163.26"
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('D29').Formula = "'136.18"
$ws.Range('E29').Value = '  +1.86%  '
$ws.Range('D30').Formula = "'5.189"
$ws.Range('E30').Value = '  +1.67%  '
$ws.Range('D31').Formula = "'7.484"
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').Value = '1.877.48'
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('D34').Formula = "'0.08598"
$ws.Range('E34').Value = '  -1.15%  '
$ws.Range('D35').Formula = "'7.141"
$ws.Range('E35').Value = '  -3.66%  '
$ws.Range('D36').Formula = "'11.64"
$ws.Range('E36').Value = '  +6.07%  '
$ws.Range('D37').Formula = "'0.2748"
$ws.Range('E37').Value = '  +2.19%  '
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').Formula = "'14.50"
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').Formula = "'0.09231"
$ws.Range('E40').Value = '  +3.27%  '
$ws.Range('D41').Formula = "'0.02734"
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('D42').Formula = "'1.468"
$ws.Range('E42').Value = '  +0.76%  '
$ws.Range('D43').Formula = "'0.7670"
$ws.Range('E43').Value = '  +1.43%  '
$ws.Range('D44').Formula = "'16.09"
$ws.Range('E44').Value = '  +4.64%  '
$ws.Range('D45').Formula = "'0.7202"
$ws.Range('E45').Value = '  +1.05%  '
$ws.Range('D46').Formula = "'2.575"
$ws.Range('D47').Formula = "'4.220"
$ws.Range('E47').Value = '  +2.05%  '
$ws.Range('D48').Formula = "'1.0000"
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('D49').Formula = "'1.324"
$ws.Range('E49').Value = '  +2.99%  '
$ws.Range('D50').Formula = "'139.42"
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('D51').Formula = "'0.07981"
$ws.Range('E51').Value = '  +0.65%  '
